$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so
# values like "44.261.16" or "37.42" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.261.16"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "2.264.41"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "320.46"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "102.50"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "37.42"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "2.599.75"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "0.865"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "2.262.97"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "44.142.86"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "13.50"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D22").Value = "65.87"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "3.16"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "236.42"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  +4.04%  "
$ws.Range("D28").Value = "39.13"
$ws.Range("E28").Value = "  +8.01%  "
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").Value = "163.07"
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("D32").Value = "20.30"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "2.68"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("E35").Value = "  +10.30%  "
$ws.Range("D36").Value = "2.00"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").Value = "3.08"
$ws.Range("E37").Value = "  -6.50%  "
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "16.71"
$ws.Range("E39").Value = "  +19.56%  "
$ws.Range("D40").Value = "3.73"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "4.24"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "1.784.87"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "82.71"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").Value = "75.21"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "105.34"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "58.63"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.68"
$ws.Range("E51").Value = "  +5.84%  "
